# annual_deaths.xlsx — "change name exposures categories + correct error in
# naming of canton Geneva"
#
#   1. Canton "Geneva" (English spelling) was mis-named; correct it to the
#      proper "Genève" everywhere it appears (rows 10 and 37, column A).
#   2. The "category" column used terse single-letter codes "O"/"U"; rename
#      them to the descriptive "Over 75" / "Under 75".
#   3. Leave the current selection on A10 (matches the saved view state,
#      with the window scrolled back to the top instead of being left
#      part-way down the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the canton name -------------------------------------------------
$ws.Range("A10").Value = "Genève"
$ws.Range("A37").Value = "Genève"

# --- 2. Rename the exposure categories --------------------------------------
# xlWhole (=1) match so we only touch cells whose entire content is the
# single letter, never a substring of another word (e.g. "Solothurn").
$ws.Range("E2:E28").Replace("O", "Over 75", 1)
$ws.Range("E29:E55").Replace("U", "Under 75", 1)

# --- 3. Reset the view / selection ------------------------------------------
$ws.Range("A10").Select()
